# Generate Report for Handoff
# The 504b6090-... file has just been handed off again; update its
# "Latest Handoff Datetime" cell (column D, row 5) on both language sheets.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D5").Value = "2016-03-08 18:40:00"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D5").Value = "2016-03-08 18:40:07"
